$wb = $excel.ActiveWorkbook

# Map of cell -> new value that must be updated on both the "展览" and
# "全部类型" worksheets (they contain duplicated listings).
$updates = @{
    "F2"  = 1886
    "F6"  = 758
    "F9"  = 4506
    "F11" = 355
    "F12" = 1282
    "F15" = 868
    "F17" = 478
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
